$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$journalText = "#1. Create a journal receiver: `n        `$ CRTJRNRCV JRNRCV(MYLIB/JRNRCV0001)`n#2. Create a journal: `n        `$ CRTJRN JRN(MYLIB/JOURNAL) JRNRCV(MYLIB/JRNRCV0001) MINENTDTA(*NONE)`n#3. Start journaling the file to the journal:`n        `$ STRJRNPF FILE(MYLIB/CUSTOMERS) JRN(MYLIB/JOURNAL)`n            IMAGES(*BOTH) OMTJRNE(*OPNCLO)`n#4. Display and dump the journal of the file:`n        `$ DSPJRN JRN(YMYLES/JOURNAL) FILE((YMYLES/CUSTOMERS)) OUTPUT(*PRINT)`n(p.s A journal can take more than 1 file at the same time)"

$ws.Range("A15").Value = "Journaling"
$ws.Range("B15").Value = "Basic Demo"
$ws.Range("C15").Value = $journalText

$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = $journalText

$ws.Rows.Item(15).RowHeight = 113.25

$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("C16").Select()
